$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 167; this shifts the existing row 167
# (and everything below it, down through the old last row 236) down
# by one row, extending the used range to A1:R237.
$ws.Rows(167).Insert()

# Populate the newly inserted row 167 with the new weekly record.
$ws.Range("A167").Value = 11
$ws.Range("B167").Value = "Vega Monumental Concepci$([char]0xF3)n"
$ws.Range("C167").Value = "B$([char]0xED)ob$([char]0xED)o"
$ws.Range("D167").Value = 44917
$ws.Range("E167").Value = 8
$ws.Range("F167").Value = 100112003
$ws.Range("G167").Value = "Ajo"
$ws.Range("H167").Value = "Chino"
$ws.Range("I167").Value = "Primera"
$ws.Range("J167").Value = 400
$ws.Range("K167").Value = 13000
$ws.Range("L167").Value = 14000
$ws.Range("M167").Value = 13500
$ws.Range("N167").Value = "$/caja 10 kilos"
$ws.Range("O167").Value = "China"
$ws.Range("P167").Value = 1350
$ws.Range("Q167").Value = 10
$ws.Range("R167").Value = "Hortaliza"
